$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("A11").Value = 131143990
$ws.Range("Q11").Value = 503447
$ws.Range("R11").Value = 7036031
$ws.Range("A12").Value = 131144272
$ws.Range("Q12").Value = 503267
$ws.Range("R12").Value = 7036125
$ws.Range("A13").Value = 131143987
$ws.Range("Q13").Value = 503418
$ws.Range("R13").Value = 7036017
$ws.Range("A14").Value = 131144287
$ws.Range("Q14").Value = 503222
$ws.Range("R14").Value = 7035916
$ws.Range("A22").Value = 131143998
$ws.Range("AC22").Value = "På flera granar."
$ws.Range("AJ22").Value = "gran"
$ws.Range("AK22").Value = "Picea abies"
$ws.Range("AM22").Value = ""
$ws.Range("AO22").Value = "Picea abies"
$ws.Range("B22").Value = 79245
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = "Garnlav"
$ws.Range("G22").Value = "Alectoria sarmentosa"
$ws.Range("H22").Value = "(Ach.) Ach."
$ws.Range("K22").Value = ""
$ws.Range("Q22").Value = 503444
$ws.Range("R22").Value = 7036006
$ws.Range("A23").Value = 131143972
$ws.Range("AC23").Value = "På rönn."
$ws.Range("B23").Value = 80379
$ws.Range("D23").Value = "LC"
$ws.Range("E23").Value = 6462
$ws.Range("F23").Value = "Stuplav"
$ws.Range("G23").Value = "Nephroma bellum"
$ws.Range("H23").Value = "(Spreng.) Tuck."
$ws.Range("K23").Value = "med apothecier"
$ws.Range("Q23").Value = 503377
$ws.Range("R23").Value = 7036016
$ws.Range("A24").Value = 131143976
$ws.Range("AC24").Value = ""
$ws.Range("AJ24").Value = "rönn"
$ws.Range("AK24").Value = "Sorbus aucuparia"
$ws.Range("AM24").Value = "Bark på levande träd"
$ws.Range("AO24").Value = "Bark on living woody plant # Sorbus aucuparia"
$ws.Range("Q24").Value = 503419
$ws.Range("R24").Value = 7036154
$ws.Range("A25").Value = 131143989
$ws.Range("AC25").Value = "Lunglav på en gran vid en rönn med lunglav."
$ws.Range("AM25").Value = "Gren på levande träd"
$ws.Range("AO25").Value = "Branch on living tree # Picea abies"
$ws.Range("B25").Value = 80350
$ws.Range("E25").Value = 6458
$ws.Range("F25").Value = "Lunglav"
$ws.Range("G25").Value = "Lobaria pulmonaria"
$ws.Range("H25").Value = "(L.) Hoffm."
$ws.Range("Q25").Value = 503448
$ws.Range("R25").Value = 7036030
$ws.Range("A42").Value = 131144282
$ws.Range("AC42").Value = "På rönn."
$ws.Range("AF42").Value = ""
$ws.Range("AJ42").Value = "rönn"
$ws.Range("AK42").Value = "Sorbus aucuparia"
$ws.Range("AM42").Value = "Bark på levande träd"
$ws.Range("AO42").Value = "Bark on living woody plant # Sorbus aucuparia"
$ws.Range("B42").Value = 80350
$ws.Range("E42").Value = 6458
$ws.Range("F42").Value = "Lunglav"
$ws.Range("G42").Value = "Lobaria pulmonaria"
$ws.Range("H42").Value = "(L.) Hoffm."
$ws.Range("J42").Value = ""
$ws.Range("L42").Value = ""
$ws.Range("M42").Value = ""
$ws.Range("Q42").Value = 503179
$ws.Range("R42").Value = 7035940
$ws.Range("A43").Value = 131144276
$ws.Range("AC43").Value = "På rönn."
$ws.Range("AF43").Value = ""
$ws.Range("AJ43").Value = "rönn"
$ws.Range("AK43").Value = "Sorbus aucuparia"
$ws.Range("AM43").Value = "Bark på levande träd"
$ws.Range("AO43").Value = "Bark on living woody plant # Sorbus aucuparia"
$ws.Range("B43").Value = 80350
$ws.Range("E43").Value = 6458
$ws.Range("F43").Value = "Lunglav"
$ws.Range("G43").Value = "Lobaria pulmonaria"
$ws.Range("H43").Value = "(L.) Hoffm."
$ws.Range("J43").Value = ""
$ws.Range("L43").Value = ""
$ws.Range("M43").Value = ""
$ws.Range("Q43").Value = 503155
$ws.Range("R43").Value = 7035984
$ws.Range("A44").Value = 131144267
$ws.Range("AC44").Value = "Ringhack, äldre, på en gran."
$ws.Range("AF44").Value = ""
$ws.Range("AJ44").Value = "gran"
$ws.Range("AK44").Value = "Picea abies"
$ws.Range("AM44").Value = ""
$ws.Range("AO44").Value = "Picea abies"
$ws.Range("B44").Value = 57884
$ws.Range("E44").Value = 100109
$ws.Range("F44").Value = "Tretåig hackspett"
$ws.Range("G44").Value = "Picoides tridactylus"
$ws.Range("H44").Value = "(Linnaeus, 1758)"
$ws.Range("J44").Value = ""
$ws.Range("L44").Value = ""
$ws.Range("M44").Value = "äldre spår"
$ws.Range("Q44").Value = 503564
$ws.Range("R44").Value = 7036045
$ws.Range("A45").Value = 131144263
$ws.Range("AC45").Value = "Ringhack, färska, på en gran."
$ws.Range("AF45").Value = ""
$ws.Range("AJ45").Value = "gran"
$ws.Range("AK45").Value = "Picea abies"
$ws.Range("AM45").Value = "Trädstam på levande träd"
$ws.Range("AO45").Value = "Stem on living tree # Picea abies"
$ws.Range("B45").Value = 57884
$ws.Range("E45").Value = 100109
$ws.Range("F45").Value = "Tretåig hackspett"
$ws.Range("G45").Value = "Picoides tridactylus"
$ws.Range("H45").Value = "(Linnaeus, 1758)"
$ws.Range("J45").Value = ""
$ws.Range("L45").Value = ""
$ws.Range("M45").Value = "färska spår"
$ws.Range("Q45").Value = 503363
$ws.Range("R45").Value = 7035960
$ws.Range("A75").Value = 131143970
$ws.Range("AJ75").Value = "rönn"
$ws.Range("AK75").Value = "Sorbus aucuparia"
$ws.Range("AM75").Value = "Bark på levande träd"
$ws.Range("AO75").Value = "Bark on living woody plant # Sorbus aucuparia"
$ws.Range("B75").Value = 80379
$ws.Range("D75").Value = "LC"
$ws.Range("E75").Value = 6462
$ws.Range("F75").Value = "Stuplav"
$ws.Range("G75").Value = "Nephroma bellum"
$ws.Range("H75").Value = "(Spreng.) Tuck."
$ws.Range("K75").Value = "med apothecier"
$ws.Range("Q75").Value = 503376
$ws.Range("R75").Value = 7035991
$ws.Range("A76").Value = 131144302
$ws.Range("AC76").Value = ""
$ws.Range("AJ76").Value = "gran"
$ws.Range("AK76").Value = "Picea abies"
$ws.Range("AM76").Value = ""
$ws.Range("AO76").Value = "Picea abies"
$ws.Range("B76").Value = 79245
$ws.Range("E76").Value = 6425
$ws.Range("F76").Value = "Garnlav"
$ws.Range("G76").Value = "Alectoria sarmentosa"
$ws.Range("H76").Value = "(Ach.) Ach."
$ws.Range("Q76").Value = 503606
$ws.Range("R76").Value = 7036064
$ws.Range("A77").Value = 131144273
$ws.Range("AC77").Value = "På rönn."
$ws.Range("B77").Value = 80350
$ws.Range("D77").Value = "NT"
$ws.Range("E77").Value = 6458
$ws.Range("F77").Value = "Lunglav"
$ws.Range("G77").Value = "Lobaria pulmonaria"
$ws.Range("H77").Value = "(L.) Hoffm."
$ws.Range("K77").Value = ""
$ws.Range("Q77").Value = 503267
$ws.Range("R77").Value = 7036145
